# Apply the table style change recorded in the commit:
#   ppt/slides/slide6.xml  (slide index 6, the "SOURCES OF FINANCE" slide)
#   <a:tableStyleId> changes from
#     {15377155-0E05-4A45-9B5F-DC8BB056307B}
#   to
#     {580B9779-EFE8-4AC9-9A43-8560F77D0699}
#
# In the PowerPoint object model a table's style GUID cannot be written
# through the Table.Style property directly (PowerPoint raises an error
# telling you to use ApplyStyle instead), so we use Table.ApplyStyle().
#
# (The rest of the underlying diff just swaps the raw contents of the two
# theme parts, ppt/theme/theme1.xml <-> ppt/theme/theme2.xml, with every
# relationship left pointing at the same part names. That is not something
# exposed by the Theme/ThemeVariant(s) COM surface here - there is no
# supported way to rewrite a master's/notes master's underlying <a:theme>
# document from script - so only the scriptable table-style edit is applied.)

$p = $ppt.ActivePresentation

$oldStyleId = "{15377155-0E05-4A45-9B5F-DC8BB056307B}"
$targetStyleId = "{580B9779-EFE8-4AC9-9A43-8560F77D0699}"
$changed = 0

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($targetStyleId)
                $changed = $changed + 1
                Write-Host "Updated table style on slide" $s.SlideIndex "shape" $sh.Name
            }
        }
    }
}

Write-Host "Tables updated:" $changed
